$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Resize / reposition the 4 surface charts.
#    All four charts grow from 5400000x2700000 EMU to 9000000x5400000 EMU
#    (708.6614173228346 pt x 425.1968503937008 pt). Charts 2 & 4 also
#    move from column L (index 11) to column R (index 17); charts 3 & 4
#    move from row 31 (index 30) to row 45 (index 44).
# ------------------------------------------------------------------
$newWidth = 708.6614173228346
$newHeight = 425.1968503937008

$co1 = $ws.ChartObjects(1)
$co1.Width = $newWidth
$co1.Height = $newHeight

$co2 = $ws.ChartObjects(2)
$co2.Left = $ws.Cells.Item(15, 18).Left
$co2.Top = $ws.Cells.Item(15, 18).Top
$co2.Width = $newWidth
$co2.Height = $newHeight

$co3 = $ws.ChartObjects(3)
$co3.Left = $ws.Cells.Item(45, 1).Left
$co3.Top = $ws.Cells.Item(45, 1).Top
$co3.Width = $newWidth
$co3.Height = $newHeight

$co4 = $ws.ChartObjects(4)
$co4.Left = $ws.Cells.Item(45, 18).Left
$co4.Top = $ws.Cells.Item(45, 18).Top
$co4.Width = $newWidth
$co4.Height = $newHeight

# ------------------------------------------------------------------
# 2) Fix the membership-function grid values in B2:W12 (column A holds
#    the unchanged x-axis values and is left untouched).
# ------------------------------------------------------------------
$rowsData = @{
    2  = @(0.47, 0.47, 0.47, 0.47, 0.47, 0.47, 0.47, 0.47, 0.47, 0.46, 0.45, 0.42, 0.38, 0.31, 0.19, 0, 0, 0, 0, 0, 0, 0)
    3  = @(0.48, 0.48, 0.48, 0.48, 0.48, 0.48, 0.48, 0.48, 0.48, 0.47, 0.46, 0.44, 0.4, 0.33, 0.2, 0.15, 0.15, 0.15, 0.15, 0.15, 0.15, 0.14)
    4  = @(0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.49, 0.48, 0.46, 0.42, 0.35, 0.28, 0.26, 0.26, 0.26, 0.26, 0.26, 0.25, 0.24)
    5  = @(0.53, 0.53, 0.53, 0.53, 0.53, 0.53, 0.53, 0.53, 0.53, 0.53, 0.51, 0.49, 0.46, 0.4, 0.38, 0.36, 0.36, 0.36, 0.36, 0.36, 0.33, 0.32)
    6  = @(0.58, 0.58, 0.58, 0.58, 0.58, 0.58, 0.58, 0.58, 0.58, 0.57, 0.5600000000000001, 0.55, 0.52, 0.48, 0.45, 0.44, 0.44, 0.44, 0.44, 0.41, 0.38, 0.37)
    7  = @(0.65, 0.65, 0.65, 0.65, 0.65, 0.65, 0.65, 0.65, 0.65, 0.65, 0.63, 0.61, 0.58, 0.54, 0.52, 0.5, 0.5, 0.5, 0.5, 0.45, 0.42, 0.4)
    8  = @(0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.76, 0.75, 0.73, 0.66, 0.59, 0.57, 0.5600000000000001, 0.5600000000000001, 0.5600000000000001, 0.52, 0.47, 0.44, 0.43)
    9  = @(0.86, 0.86, 0.86, 0.86, 0.86, 0.86, 0.86, 0.86, 0.86, 0.86, 0.85, 0.8, 0.74, 0.67, 0.63, 0.61, 0.61, 0.59, 0.54, 0.49, 0.46, 0.44)
    10 = @(0.93, 0.93, 0.93, 0.93, 0.93, 0.93, 0.93, 0.93, 0.93, 0.93, 0.91, 0.86, 0.8, 0.73, 0.7, 0.67, 0.67, 0.6, 0.55, 0.5, 0.47, 0.46)
    11 = @(0.98, 0.98, 0.98, 0.98, 0.98, 0.98, 0.98, 0.98, 0.98, 0.97, 0.9399999999999999, 0.89, 0.83, 0.77, 0.74, 0.71, 0.7, 0.61, 0.55, 0.5, 0.47, 0.46)
    12 = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 0.99, 0.96, 0.91, 0.85, 0.79, 0.76, 0.73, 0.7, 0.62, 0.55, 0.51, 0.48, 0.47)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $col = 2
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
